$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "70.440.38"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  +6.28%  "
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "3.645.29"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  +18.67%  "
$ws.Range("E4").Value = "  -0.19%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "595.86"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +3.56%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "185.44"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +8.70%  "
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "3.644.32"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  +18.75%  "
$ws.Range("E8").Value = "  -0.06%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.535"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +4.97%  "
$ws.Range("E10").Value = "  +9.64%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "6.54"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  +4.33%  "
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.497"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  +6.02%  "
$ws.Range("B13").Value = "ShibaInu"
$ws.Range("C13").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "0.0000256"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  +7.38%  "
$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "39.36"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  +10.25%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "4.245.44"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +18.45%  "
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "3.638.64"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  +18.34%  "
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "70.231.73"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  +6.00%  "
$ws.Range("E18").Value = "  +2.14%  "
$ws.Range("E19").Value = "  +8.30%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "17.35"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +4.45%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "511.15"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +5.43%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "9.22"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +20.64%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "0.751"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +9.74%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "88.34"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +7.35%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "13.56"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +7.61%  "
$ws.Range("E26").Value = "  +9.23%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "10.78"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  +7.42%  "
$ws.Range("E28").Value = "  +0.09%  "
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "2.56"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  +14.07%  "
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "8.24"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +5.21%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "32.03"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  +15.93%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "2.76"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  +6.57%  "
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "0.0000109"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  +18.90%  "
$ws.Range("E34").Value = "  +6.00%  "
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "0.997"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -0.32%  "
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "6.16"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  +11.01%  "
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "1.02"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +8.66%  "
$ws.Range("E38").Value = "  +11.59%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "2.12"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  +8.30%  "
$ws.Range("B40").Value = "Arweave"
$ws.Range("C40").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "47.14"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -1.26%  "
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "50.86"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  +3.76%  "
$ws.Range("E42").Value = "  +5.46%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "8.92"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +8.55%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "3.146.33"
$cell.Style = "Normal"
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "2.80"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  +10.44%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "406.77"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  +11.81%  "
$ws.Range("E47").Value = "  +6.66%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "27.96"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  +15.48%  "
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "134.65"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +0.00%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "2.47"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  +15.42%  "
